$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, [string]$text)
    $scratch = $ws.Range("A200")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.EntireRow.Delete()
}

$ws.Range("D2").Value = "25.984.99"
$ws.Range("E2").Value = "  +5.53%  "
$ws.Range("D3").Value = "1.716.51"
$ws.Range("E3").Value = "  +3.49%  "
Set-TextValue $ws "D4" "1.002"
$ws.Range("E4").Value = "  +0.22%  "
Set-TextValue $ws "D5" "330.86"
$ws.Range("E5").Value = "  +2.94%  "
Set-TextValue $ws "D6" "1.000"
$ws.Range("E6").Value = "  +0.24%  "
Set-TextValue $ws "D7" "0.3693"
$ws.Range("E7").Value = "  +1.12%  "
Set-TextValue $ws "D8" "49.62"
$ws.Range("E8").Value = "  +6.06%  "
Set-TextValue $ws "D9" "0.3324"
$ws.Range("E9").Value = "  +1.48%  "
Set-TextValue $ws "D10" "1.182"
$ws.Range("E10").Value = "  +4.29%  "
Set-TextValue $ws "D11" "0.07481"
$ws.Range("E11").Value = "  +5.90%  "
Set-TextValue $ws "D12" "1.000"
$ws.Range("E12").Value = "  +0.36%  "
Set-TextValue $ws "D13" "6.242"
$ws.Range("E13").Value = "  +4.05%  "
Set-TextValue $ws "D14" "20.10"
$ws.Range("E14").Value = "  +2.76%  "
Set-TextValue $ws "D15" "6.905"
$ws.Range("E15").Value = "  +4.15%  "
$ws.Range("D16").Value = "1.715.19"
$ws.Range("E16").Value = "  +3.67%  "
Set-TextValue $ws "D17" "0.00001077"
$ws.Range("E17").Value = "  +2.80%  "
Set-TextValue $ws "D18" "0.06636"
$ws.Range("E18").Value = "  +0.11%  "
Set-TextValue $ws "D19" "82.02"
$ws.Range("E19").Value = "  +3.58%  "
Set-TextValue $ws "D20" "0.9994"
$ws.Range("E20").Value = "  +0.23%  "
Set-TextValue $ws "D21" "16.32"
$ws.Range("E21").Value = "  +3.05%  "
Set-TextValue $ws "D22" "6.078"
$ws.Range("E22").Value = "  +1.88%  "
Set-TextValue $ws "D23" "13.04"
$ws.Range("E23").Value = "  +3.16%  "
$ws.Range("D24").Value = "25.924.18"
$ws.Range("E24").Value = "  +5.43%  "
Set-TextValue $ws "D25" "2.473"
$ws.Range("E25").Value = "  +0.28%  "
Set-TextValue $ws "D26" "2.484"
$ws.Range("E26").Value = "  +4.09%  "
Set-TextValue $ws "D27" "150.01"
$ws.Range("E27").Value = "  +1.24%  "
Set-TextValue $ws "D28" "19.27"
$ws.Range("E28").Value = "  +3.10%  "
Set-TextValue $ws "D29" "1.302"
$ws.Range("E29").Value = "  +6.64%  "
$ws.Range("D30").Value = "1.908.92"
$ws.Range("E30").Value = "  +3.84%  "
Set-TextValue $ws "D31" "128.81"
$ws.Range("E31").Value = "  +2.95%  "
Set-TextValue $ws "D32" "4.107"
$ws.Range("E32").Value = "  +1.11%  "
Set-TextValue $ws "D33" "5.956"
$ws.Range("E33").Value = "  +1.66%  "
Set-TextValue $ws "D34" "0.08539"
$ws.Range("E34").Value = "  +0.54%  "
Set-TextValue $ws "D35" "1.719"
$ws.Range("E35").Value = "  +2.39%  "
Set-TextValue $ws "D36" "12.91"
$ws.Range("E36").Value = "  +4.57%  "
Set-TextValue $ws "D37" "5.352"
$ws.Range("E37").Value = "  +2.16%  "
Set-TextValue $ws "D38" "1.281"
$ws.Range("E38").Value = "  +0.35%  "
Set-TextValue $ws "D39" "0.06221"
$ws.Range("E39").Value = "  +2.76%  "
$ws.Range("E40").Value = "  +2.27%  "
Set-TextValue $ws "D41" "0.2139"
$ws.Range("E41").Value = "  +2.97%  "
Set-TextValue $ws "D42" "8.544"
$ws.Range("E42").Value = "  +4.58%  "
Set-TextValue $ws "D43" "14.45"
$ws.Range("E43").Value = "  +13.53%  "
Set-TextValue $ws "D44" "0.6160"
$ws.Range("E44").Value = "  +3.77%  "
Set-TextValue $ws "D45" "0.9998"
$ws.Range("E45").Value = "  +0.29%  "
Set-TextValue $ws "D46" "3.842"
$ws.Range("E46").Value = "  -0.65%  "
Set-TextValue $ws "D47" "0.5865"
$ws.Range("E47").Value = "  +3.93%  "
Set-TextValue $ws "D48" "127.16"
$ws.Range("E48").Value = "  +2.12%  "
Set-TextValue $ws "D49" "2.009"
$ws.Range("E49").Value = "  +2.38%  "
Set-TextValue $ws "D50" "0.07250"
$ws.Range("E50").Value = "  +3.90%  "
Set-TextValue $ws "D51" "77.05"
$ws.Range("E51").Value = "  +2.84%  "
